$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.106.57'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.97%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.018.71'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -1.61%  '
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.21'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.70%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.609'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.63%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '55.78'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.92%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.374'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -3.23%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0779'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -3.20%  '
$ws.Range("E11").Value = '  -4.00%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.315.63'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.58%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.13'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -2.85%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.90'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -3.49%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.739'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -2.11%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.18'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.26%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.009.32'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.49%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.067.06'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.73%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.24'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +2.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '68.97'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.17%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0813'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -4.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '223.15'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.17%  '
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("E24").Value = '  +2.44%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.19'
$ws.Range("D25").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.03'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.86%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.01'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -5.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.127'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.37%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.60'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.66%  '
$ws.Range("E30").Value = '  -4.56%  '
$ws.Range("E31").Value = '  -0.90%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.42'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.56%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0601'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.79%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.46'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -1.36%  '
$ws.Range("B35").Value = 'LidoDAOToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.35'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.98%  '
$ws.Range("B36").Value = 'WEMIXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.87'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +2.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.15'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.36%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.43'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.13%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.465.64'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.98%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0213'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -3.72%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '94.56'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.46%  '
$ws.Range("E43").Value = '  -3.04%  '
$ws.Range("E44").Value = '  -3.32%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.01'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -4.68%  '
$ws.Range("E46").Value = '  +7.77%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.12'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.50%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.47%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.06'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.17%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.92'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.40%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.201.93'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.67%  '
